$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix spelling: "Category Definitino" -> "Category Definition"
$ws.Range("D2").Value = "Category Definition"

# Reset the view so the top-left visible cell is A1 and the active
# selection lands on H4 (matches the saved view state after the edit).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("H4").Select()
